$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "835×5="
$t.Cell(1,2).Range.Text  = "614×4="
$t.Cell(1,3).Range.Text  = "931×3="
$t.Cell(1,4).Range.Text  = "556×4="
$t.Cell(1,5).Range.Text  = "350×3="

$t.Cell(5,1).Range.Text  = "966×3="
$t.Cell(5,2).Range.Text  = "396×5="
$t.Cell(5,3).Range.Text  = "136×4="
$t.Cell(5,4).Range.Text  = "754×3="
$t.Cell(5,5).Range.Text  = "217×9="

$t.Cell(10,1).Range.Text = "257×5="
$t.Cell(10,2).Range.Text = "233×7="
$t.Cell(10,3).Range.Text = "265×6="
$t.Cell(10,4).Range.Text = "690×7="
$t.Cell(10,5).Range.Text = "384×6="

$t.Cell(15,1).Range.Text = "160×9="
$t.Cell(15,2).Range.Text = "733×2="
$t.Cell(15,3).Range.Text = "972×2="
$t.Cell(15,4).Range.Text = "330×6="
$t.Cell(15,5).Range.Text = "943×4="

$t.Cell(20,1).Range.Text = "866×5="
$t.Cell(20,2).Range.Text = "152×2="
$t.Cell(20,3).Range.Text = "153×6="
$t.Cell(20,4).Range.Text = "883×4="
$t.Cell(20,5).Range.Text = "335×7="
